$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 231
$ws1.Range("F3").Value = 55045
$ws1.Range("F4").Value = 1274
$ws1.Range("F6").Value = 874
$ws1.Range("F8").Value = 1153
$ws1.Range("F9").Value = 1433
$ws1.Range("F10").Value = 140
$ws1.Range("F11").Value = 47
$ws1.Range("F12").Value = 248
$ws1.Range("F13").Value = 413
$ws1.Range("F14").Value = 80
$ws1.Range("F15").Value = 43
$ws1.Range("F17").Value = 78
$ws1.Range("F18").Value = 73
$ws1.Range("F19").Value = 5684
$ws1.Range("F20").Value = 43
$ws1.Range("F21").Value = 5515
$ws1.Range("F22").Value = 9563
$ws1.Range("F25").Value = 163
$ws1.Range("F26").Value = 252
$ws1.Range("F27").Value = 472
$ws1.Range("F28").Value = 149
$ws1.Range("F29").Value = 118
$ws1.Range("F30").Value = 4320
$ws1.Range("F31").Value = 333

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 60
$ws2.Range("F20").Value = 10

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 231
$ws4.Range("F5").Value = 1274
$ws4.Range("F8").Value = 874
$ws4.Range("F10").Value = 1153
$ws4.Range("F12").Value = 1433
$ws4.Range("F14").Value = 140
$ws4.Range("F15").Value = 248
$ws4.Range("F17").Value = 413
$ws4.Range("F18").Value = 80
$ws4.Range("F19").Value = 43
$ws4.Range("F22").Value = 78
$ws4.Range("F23").Value = 73
$ws4.Range("F24").Value = 5684
$ws4.Range("F25").Value = 43
$ws4.Range("F26").Value = 5515
$ws4.Range("F27").Value = 9563
$ws4.Range("F31").Value = 163
$ws4.Range("F32").Value = 252
$ws4.Range("F34").Value = 472
$ws4.Range("F37").Value = 149
$ws4.Range("F38").Value = 118
$ws4.Range("F39").Value = 4320
$ws4.Range("F40").Value = 60
$ws4.Range("F46").Value = 333
$ws4.Range("F47").Value = 10
